# The edit reshuffles the data rows of the sheet (rows 2-33): the values
# that used to live in one row now live in another row, while row 1
# (headers) is untouched. Only the columns whose content actually differs
# from row to row need to move; columns that hold the same constant value
# in every data row are left completely alone so the round-trip through
# the Excel object model can't introduce incidental side effects (e.g.
# date-like text being re-typed as a real date) in columns the diff never
# touches.
#
# destination row-offset (1 == row 2 ... 32 == row 33) -> source row-offset
$map = @{1=25; 2=13; 3=12; 4=17; 5=20; 6=29; 7=19; 8=32; 9=21; 10=3; 11=30; 12=31; 13=14; 14=1; 15=4; 16=7; 17=15; 18=6; 19=24; 20=2; 21=11; 22=23; 23=16; 24=27; 25=10; 26=8; 27=26; 28=5; 29=9; 30=28; 31=18; 32=22}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowCount = 32
$firstRow = 2
$lastRow = 33

# Columns whose value varies between data rows (everything else is a
# constant repeated on every row, so it doesn't need to be touched).
$cols = @("A","B","D","E","F","G","H","I","M","N","Q","R","AC")

foreach ($col in $cols) {
    $rangeAddr = "$col$($firstRow):$col$($lastRow)"
    $src = $ws.Range($rangeAddr).Value()
    $dst = $ws.Range($rangeAddr).Value()

    for ($r = 1; $r -le $rowCount; $r++) {
        $sr = $map[$r]
        $dst[$r, 1] = $src[$sr, 1]
    }

    $ws.Range($rangeAddr).Value = $dst
}
